$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Subpopulation" (sub_pop_section*) and "LOT" (lot_section*) page-name
# values in the AddtionalParam/sectionname columns are being renamed to the
# new "Population Filter 1" / "Population Filter 2" naming scheme.

# Row 2: sub_pop_section1 / sub_pop_section1_checkbox / sub_pop_section
#   -> pop_filter1_section1 / pop_filter1_section1_checkbox / pop_filter1_section
$ws.Range("F2").Value = "pop_filter1_section1"
$ws.Range("G2").Value = "pop_filter1_section1_checkbox"
$ws.Range("H2").Value = "pop_filter1_section"

# Row 3: sub_pop_section2 / sub_pop_section2_checkbox / sub_pop_section
#   -> pop_filter1_section2 / pop_filter1_section2_checkbox / pop_filter1_section
$ws.Range("F3").Value = "pop_filter1_section2"
$ws.Range("G3").Value = "pop_filter1_section2_checkbox"
$ws.Range("H3").Value = "pop_filter1_section"

# Row 4: lot_section2 / lot_section2_checkbox / lot_section
#   -> pop_filter2_section2 / pop_filter2_section2_checkbox / pop_filter2_section
$ws.Range("F4").Value = "pop_filter2_section2"
$ws.Range("G4").Value = "pop_filter2_section2_checkbox"
$ws.Range("H4").Value = "pop_filter2_section"

# Update the sheet's selection/view state to match the saved workbook state:
# active cell F2, selected range F2:H4, scrolled so column C is at the
# left edge of the window.
$ws.Range("F2:H4").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
